# Update "F" column ("想去人数" / interest counts) across the four sheets
# (展览, 演出, 本地生活, 全部类型) to reflect the regenerated gh-pages data
# snapshot at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 3308
$ws.Range("F5").Value = 338
$ws.Range("F6").Value = 7729
$ws.Range("F8").Value = 725
$ws.Range("F9").Value = 1137
$ws.Range("F10").Value = 1060
$ws.Range("F12").Value = 13
$ws.Range("F14").Value = 1750
$ws.Range("F16").Value = 6178
$ws.Range("F19").Value = 1008
$ws.Range("F23").Value = 6314
$ws.Range("F24").Value = 6466
$ws.Range("F31").Value = 1065
$ws.Range("F32").Value = 1037
$ws.Range("F33").Value = 109
$ws.Range("F34").Value = 109
$ws.Range("F44").Value = 712
$ws.Range("F47").Value = 3233

$ws = $wb.Worksheets.Item(2)
$ws.Range("F10").Value = 283
$ws.Range("F17").Value = 36
$ws.Range("F24").Value = 6577

$ws = $wb.Worksheets.Item(3)
$ws.Range("F8").Value = 2140
$ws.Range("F9").Value = 8909
$ws.Range("F11").Value = 81

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 3308
$ws.Range("F4").Value = 338
$ws.Range("F5").Value = 7729
$ws.Range("F10").Value = 81
$ws.Range("F11").Value = 1137
$ws.Range("F12").Value = 1060
$ws.Range("F16").Value = 283
$ws.Range("F22").Value = 6314
$ws.Range("F23").Value = 6466
$ws.Range("F30").Value = 1065
$ws.Range("F31").Value = 109
$ws.Range("F32").Value = 109
$ws.Range("F43").Value = 712
$ws.Range("F45").Value = 3233

